$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -4
$ws.Range("F6").Value = -3
$ws.Range("F12").Value = -4
$ws.Range("F13").Value = 4
$ws.Range("F14").Value = -5
